$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# New demo rows, appended to the Tableau1 table (id, type, name, path, last_update)
$newRows = @(
    @("bevnat_info",     "pdf", "BEVNAT: Fiche signalétique",       "data/pdf/bevnat_info.pdf",     1706212962),
    @("statpop_info",    "pdf", "STATPOP: Fiche signalétique",      "data/pdf/statpop_info.pdf",    1606212962),
    @("bevnat_variable", "pdf", "BEVNAT: Liste des variables",      "data/pdf/bevnat_variable.pdf", 1606212963),
    @("pop_com_1",       "pdf", "Communiqué de presse population",  "data/pdf/pop_com_1.pdf",       1724323867)
)

foreach ($r in $newRows) {
    $row = $tbl.ListRows.Add()
    $row.Range.Cells.Item(1, 1).Value = $r[0]
    $row.Range.Cells.Item(1, 2).Value = $r[1]
    $row.Range.Cells.Item(1, 3).Value = $r[2]
    $row.Range.Cells.Item(1, 4).Value = $r[3]
    $row.Range.Cells.Item(1, 5).Value = $r[4]
}

# Widen the "name" column (C) to fit the new, longer labels
$ws.Columns.Item(3).ColumnWidth = 21.2

# Move the active selection the way the author left it
$ws.Range("E11").Select()
